$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dnajb11"
$ws.Cells.Item(2, 3).Value = "Prtg"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 10.70240133333333
$ws.Cells.Item(2, 8).Value = 32.107204
$ws.Cells.Item(2, 9).Value = 0.1749841631264304
$ws.Cells.Item(2, 10).Value = 0.1749841631264304
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.1033963333333333
$ws.Cells.Item(2, 14).Value = 0.310189
$ws.Cells.Item(2, 15).Value = 0.0818371535912535
$ws.Cells.Item(2, 16).Value = 0.08183715359125351
$ws.Cells.Item(2, 17).Value = 1.106589055728445
$ws.Cells.Item(2, 18).Value = 9.959301501556
$ws.Cells.Item(2, 19).Value = 0.01432020583381464
$ws.Cells.Item(2, 20).Value = 0.01432020583381465

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Dnajb11"
$ws.Cells.Item(3, 3).Value = "Prtg"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 10.70240133333333
$ws.Cells.Item(3, 8).Value = 32.107204
$ws.Cells.Item(3, 9).Value = 0.1749841631264304
$ws.Cells.Item(3, 10).Value = 0.1749841631264304
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.7274516666666667
$ws.Cells.Item(3, 14).Value = 2.182355
$ws.Cells.Item(3, 15).Value = 0.5757706473332067
$ws.Cells.Item(3, 16).Value = 0.5757706473332068
$ws.Cells.Item(3, 17).Value = 7.78547968726889
$ws.Cells.Item(3, 18).Value = 70.06931718542002
$ws.Cells.Item(3, 19).Value = 0.1007507448763643
$ws.Cells.Item(3, 20).Value = 0.1007507448763643

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Dnajb11"
$ws.Cells.Item(4, 3).Value = "Prtg"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 10.70240133333333
$ws.Cells.Item(4, 8).Value = 32.107204
$ws.Cells.Item(4, 9).Value = 0.1749841631264304
$ws.Cells.Item(4, 10).Value = 0.1749841631264304
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.000222
$ws.Cells.Item(4, 14).Value = 0.000666
$ws.Cells.Item(4, 15).Value = 0.0001757107579307288
$ws.Cells.Item(4, 16).Value = 0.0001757107579307288
$ws.Cells.Item(4, 17).Value = 0.002375933096
$ws.Cells.Item(4, 18).Value = 0.021383397864
$ws.Cells.Item(4, 19).Value = 0.00003074659992881938
$ws.Cells.Item(4, 20).Value = 0.00003074659992881938

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Dnajb11"
$ws.Cells.Item(5, 3).Value = "Prtg"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 10.70240133333333
$ws.Cells.Item(5, 8).Value = 32.107204
$ws.Cells.Item(5, 9).Value = 0.1749841631264304
$ws.Cells.Item(5, 10).Value = 0.1749841631264304
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.43237
$ws.Cells.Item(5, 14).Value = 1.29711
$ws.Cells.Item(5, 15).Value = 0.342216488317609
$ws.Cells.Item(5, 16).Value = 0.3422164883176091
$ws.Cells.Item(5, 17).Value = 4.627397264493333
$ws.Cells.Item(5, 18).Value = 41.64657538044001
$ws.Cells.Item(5, 19).Value = 0.05988246581632267
$ws.Cells.Item(5, 20).Value = 0.05988246581632268

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Dnajb11"
$ws.Cells.Item(6, 3).Value = "Prtg"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 18.645164
$ws.Cells.Item(6, 8).Value = 55.935492
$ws.Cells.Item(6, 9).Value = 0.3048482594960664
$ws.Cells.Item(6, 10).Value = 0.3048482594960664
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.1033963333333333
$ws.Cells.Item(6, 14).Value = 0.310189
$ws.Cells.Item(6, 15).Value = 0.0818371535912535
$ws.Cells.Item(6, 16).Value = 0.08183715359125351
$ws.Cells.Item(6, 17).Value = 1.927841591998666
$ws.Cells.Item(6, 18).Value = 17.350574327988
$ws.Cells.Item(6, 19).Value = 0.02494791383440589
$ws.Cells.Item(6, 20).Value = 0.02494791383440589

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Dnajb11"
$ws.Cells.Item(7, 3).Value = "Prtg"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 18.645164
$ws.Cells.Item(7, 8).Value = 55.935492
$ws.Cells.Item(7, 9).Value = 0.3048482594960664
$ws.Cells.Item(7, 10).Value = 0.3048482594960664
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.7274516666666667
$ws.Cells.Item(7, 14).Value = 2.182355
$ws.Cells.Item(7, 15).Value = 0.5757706473332067
$ws.Cells.Item(7, 16).Value = 0.5757706473332068
$ws.Cells.Item(7, 17).Value = 13.56345562707333
$ws.Cells.Item(7, 18).Value = 122.07110064366
$ws.Cells.Item(7, 19).Value = 0.1755226797084515
$ws.Cells.Item(7, 20).Value = 0.1755226797084516

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Dnajb11"
$ws.Cells.Item(8, 3).Value = "Prtg"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 18.645164
$ws.Cells.Item(8, 8).Value = 55.935492
$ws.Cells.Item(8, 9).Value = 0.3048482594960664
$ws.Cells.Item(8, 10).Value = 0.3048482594960664
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.000222
$ws.Cells.Item(8, 14).Value = 0.000666
$ws.Cells.Item(8, 15).Value = 0.0001757107579307288
$ws.Cells.Item(8, 16).Value = 0.0001757107579307288
$ws.Cells.Item(8, 17).Value = 0.004139226407999999
$ws.Cells.Item(8, 18).Value = 0.037253037672
$ws.Cells.Item(8, 19).Value = 0.00005356511872991732
$ws.Cells.Item(8, 20).Value = 0.00005356511872991733

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Dnajb11"
$ws.Cells.Item(9, 3).Value = "Prtg"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 18.645164
$ws.Cells.Item(9, 8).Value = 55.935492
$ws.Cells.Item(9, 9).Value = 0.3048482594960664
$ws.Cells.Item(9, 10).Value = 0.3048482594960664
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.43237
$ws.Cells.Item(9, 14).Value = 1.29711
$ws.Cells.Item(9, 15).Value = 0.342216488317609
$ws.Cells.Item(9, 16).Value = 0.3422164883176091
$ws.Cells.Item(9, 17).Value = 8.061609558679999
$ws.Cells.Item(9, 18).Value = 72.55448602812
$ws.Cells.Item(9, 19).Value = 0.1043241008344791
$ws.Cells.Item(9, 20).Value = 0.1043241008344791

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Dnajb11"
$ws.Cells.Item(10, 3).Value = "Prtg"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 14.50031533333333
$ws.Cells.Item(10, 8).Value = 43.500946
$ws.Cells.Item(10, 9).Value = 0.2370800220105756
$ws.Cells.Item(10, 10).Value = 0.2370800220105756
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.1033963333333333
$ws.Cells.Item(10, 14).Value = 0.310189
$ws.Cells.Item(10, 15).Value = 0.0818371535912535
$ws.Cells.Item(10, 16).Value = 0.08183715359125351
$ws.Cells.Item(10, 17).Value = 1.499279437643778
$ws.Cells.Item(10, 18).Value = 13.493514938794
$ws.Cells.Item(10, 19).Value = 0.01940195417469723
$ws.Cells.Item(10, 20).Value = 0.01940195417469724

# Row 11
$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value = "Dnajb11"
$ws.Cells.Item(11, 3).Value = "Prtg"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 14.50031533333333
$ws.Cells.Item(11, 8).Value = 43.500946
$ws.Cells.Item(11, 9).Value = 0.2370800220105756
$ws.Cells.Item(11, 10).Value = 0.2370800220105756
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.7274516666666667
$ws.Cells.Item(11, 14).Value = 2.182355
$ws.Cells.Item(11, 15).Value = 0.5757706473332067
$ws.Cells.Item(11, 16).Value = 0.5757706473332068
$ws.Cells.Item(11, 17).Value = 10.54827855642556
$ws.Cells.Item(11, 18).Value = 94.93450700783001
$ws.Cells.Item(11, 19).Value = 0.1365037177428
$ws.Cells.Item(11, 20).Value = 0.1365037177428

# Row 12
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Dnajb11"
$ws.Cells.Item(12, 3).Value = "Prtg"
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 14.50031533333333
$ws.Cells.Item(12, 8).Value = 43.500946
$ws.Cells.Item(12, 9).Value = 0.2370800220105756
$ws.Cells.Item(12, 10).Value = 0.2370800220105756
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.000222
$ws.Cells.Item(12, 14).Value = 0.000666
$ws.Cells.Item(12, 15).Value = 0.0001757107579307288
$ws.Cells.Item(12, 16).Value = 0.0001757107579307288
$ws.Cells.Item(12, 17).Value = 0.003219070004
$ws.Cells.Item(12, 18).Value = 0.028971630036
$ws.Cells.Item(12, 19).Value = 0.0000416575103577121
$ws.Cells.Item(12, 20).Value = 0.0000416575103577121

# Row 13
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Dnajb11"
$ws.Cells.Item(13, 3).Value = "Prtg"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 14.50031533333333
$ws.Cells.Item(13, 8).Value = 43.500946
$ws.Cells.Item(13, 9).Value = 0.2370800220105756
$ws.Cells.Item(13, 10).Value = 0.2370800220105756
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.43237
$ws.Cells.Item(13, 14).Value = 1.29711
$ws.Cells.Item(13, 15).Value = 0.342216488317609
$ws.Cells.Item(13, 16).Value = 0.3422164883176091
$ws.Cells.Item(13, 17).Value = 6.269501340673333
$ws.Cells.Item(13, 18).Value = 56.42551206606
$ws.Cells.Item(13, 19).Value = 0.08113269258272063
$ws.Cells.Item(13, 20).Value = 0.08113269258272064

# Row 14
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Dnajb11"
$ws.Cells.Item(14, 3).Value = "Prtg"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 5.766972666666667
$ws.Cells.Item(14, 8).Value = 17.300918
$ws.Cells.Item(14, 9).Value = 0.09428994992989723
$ws.Cells.Item(14, 10).Value = 0.09428994992989723
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.1033963333333333
$ws.Cells.Item(14, 14).Value = 0.310189
$ws.Cells.Item(14, 15).Value = 0.0818371535912535
$ws.Cells.Item(14, 16).Value = 0.08183715359125351
$ws.Cells.Item(14, 17).Value = 0.5962838281668889
$ws.Cells.Item(14, 18).Value = 5.366554453501999
$ws.Cells.Item(14, 19).Value = 0.007716421114524602
$ws.Cells.Item(14, 20).Value = 0.007716421114524603

# Row 15
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Dnajb11"
$ws.Cells.Item(15, 3).Value = "Prtg"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 5.766972666666667
$ws.Cells.Item(15, 8).Value = 17.300918
$ws.Cells.Item(15, 9).Value = 0.09428994992989723
$ws.Cells.Item(15, 10).Value = 0.09428994992989723
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.7274516666666667
$ws.Cells.Item(15, 14).Value = 2.182355
$ws.Cells.Item(15, 15).Value = 0.5757706473332067
$ws.Cells.Item(15, 16).Value = 0.5757706473332068
$ws.Cells.Item(15, 17).Value = 4.195193877987778
$ws.Cells.Item(15, 18).Value = 37.75674490189
$ws.Cells.Item(15, 19).Value = 0.05428938550815258
$ws.Cells.Item(15, 20).Value = 0.05428938550815259

# Row 16
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Dnajb11"
$ws.Cells.Item(16, 3).Value = "Prtg"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 5.766972666666667
$ws.Cells.Item(16, 8).Value = 17.300918
$ws.Cells.Item(16, 9).Value = 0.09428994992989723
$ws.Cells.Item(16, 10).Value = 0.09428994992989723
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.000222
$ws.Cells.Item(16, 14).Value = 0.000666
$ws.Cells.Item(16, 15).Value = 0.0001757107579307288
$ws.Cells.Item(16, 16).Value = 0.0001757107579307288
$ws.Cells.Item(16, 17).Value = 0.001280267932
$ws.Cells.Item(16, 18).Value = 0.011522411388
$ws.Cells.Item(16, 19).Value = 0.00001656775856743271
$ws.Cells.Item(16, 20).Value = 0.00001656775856743271

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Dnajb11"
$ws.Cells.Item(17, 3).Value = "Prtg"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 5.766972666666667
$ws.Cells.Item(17, 8).Value = 17.300918
$ws.Cells.Item(17, 9).Value = 0.09428994992989723
$ws.Cells.Item(17, 10).Value = 0.09428994992989723
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.43237
$ws.Cells.Item(17, 14).Value = 1.29711
$ws.Cells.Item(17, 15).Value = 0.342216488317609
$ws.Cells.Item(17, 16).Value = 0.3422164883176091
$ws.Cells.Item(17, 17).Value = 2.493465971886667
$ws.Cells.Item(17, 18).Value = 22.44119374698
$ws.Cells.Item(17, 19).Value = 0.03226757554865262
$ws.Cells.Item(17, 20).Value = 0.03226757554865262

# Row 18
$ws.Cells.Item(18, 1).Value = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value = "Dnajb11"
$ws.Cells.Item(18, 3).Value = "Prtg"
$ws.Cells.Item(18, 4).Value = "ECs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 11.54726066666667
$ws.Cells.Item(18, 8).Value = 34.641782
$ws.Cells.Item(18, 9).Value = 0.1887976054370303
$ws.Cells.Item(18, 10).Value = 0.1887976054370303
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.1033963333333333
$ws.Cells.Item(18, 14).Value = 0.310189
$ws.Cells.Item(18, 15).Value = 0.0818371535912535
$ws.Cells.Item(18, 16).Value = 0.08183715359125351
$ws.Cells.Item(18, 17).Value = 1.193944412977555
$ws.Cells.Item(18, 18).Value = 10.745499716798
$ws.Cells.Item(18, 19).Value = 0.01545065863381112
$ws.Cells.Item(18, 20).Value = 0.01545065863381113

# Row 19
$ws.Cells.Item(19, 1).Value = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value = "Dnajb11"
$ws.Cells.Item(19, 3).Value = "Prtg"
$ws.Cells.Item(19, 4).Value = "FAPs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 11.54726066666667
$ws.Cells.Item(19, 8).Value = 34.641782
$ws.Cells.Item(19, 9).Value = 0.1887976054370303
$ws.Cells.Item(19, 10).Value = 0.1887976054370303
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.7274516666666667
$ws.Cells.Item(19, 14).Value = 2.182355
$ws.Cells.Item(19, 15).Value = 0.5757706473332067
$ws.Cells.Item(19, 16).Value = 0.5757706473332068
$ws.Cells.Item(19, 17).Value = 8.400074017401112
$ws.Cells.Item(19, 18).Value = 75.60066615661
$ws.Cells.Item(19, 19).Value = 0.1087041194974383
$ws.Cells.Item(19, 20).Value = 0.1087041194974383

# Row 20
$ws.Cells.Item(20, 1).Value = "Resolving-Mac"
$ws.Cells.Item(20, 2).Value = "Dnajb11"
$ws.Cells.Item(20, 3).Value = "Prtg"
$ws.Cells.Item(20, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 11.54726066666667
$ws.Cells.Item(20, 8).Value = 34.641782
$ws.Cells.Item(20, 9).Value = 0.1887976054370303
$ws.Cells.Item(20, 10).Value = 0.1887976054370303
$ws.Cells.Item(20, 11).Value = 1
$ws.Cells.Item(20, 12).Value = 0.3333333333333333
$ws.Cells.Item(20, 13).Value = 0.000222
$ws.Cells.Item(20, 14).Value = 0.000666
$ws.Cells.Item(20, 15).Value = 0.0001757107579307288
$ws.Cells.Item(20, 16).Value = 0.0001757107579307288
$ws.Cells.Item(20, 17).Value = 0.002563491868
$ws.Cells.Item(20, 18).Value = 0.023071426812
$ws.Cells.Item(20, 19).Value = 0.00003317377034684727
$ws.Cells.Item(20, 20).Value = 0.00003317377034684728

# Row 21
$ws.Cells.Item(21, 1).Value = "Resolving-Mac"
$ws.Cells.Item(21, 2).Value = "Dnajb11"
$ws.Cells.Item(21, 3).Value = "Prtg"
$ws.Cells.Item(21, 4).Value = "MuSCs"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 11.54726066666667
$ws.Cells.Item(21, 8).Value = 34.641782
$ws.Cells.Item(21, 9).Value = 0.1887976054370303
$ws.Cells.Item(21, 10).Value = 0.1887976054370303
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 0.43237
$ws.Cells.Item(21, 14).Value = 1.29711
$ws.Cells.Item(21, 15).Value = 0.342216488317609
$ws.Cells.Item(21, 16).Value = 0.3422164883176091
$ws.Cells.Item(21, 17).Value = 4.992689094446666
$ws.Cells.Item(21, 18).Value = 44.93420185002
$ws.Cells.Item(21, 19).Value = 0.06460965353543403
$ws.Cells.Item(21, 20).Value = 0.06460965353543405
